# Update the "取得日時" (acquisition timestamp) column on the ランサーズ sheet
# for rows 2-10 from 2025-11-25 06:28:40 to 2025-11-25 06:37:10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-25 06:37:10"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
